# Append a new scrape result (2025-09-16 01:14:16) to the top of the
# "ランサーズ" listing sheet, pushing all previously-scraped rows down by
# one, refreshing their "取得日時" (captured-at) timestamp, and widening
# column D slightly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# 1. Insert a brand-new row above the current row 2 (the first data row).
#    This shifts the existing data rows (2-10) down to rows 3-11 while
#    preserving their cell contents/formatting.
$ws.Rows.Item(2).Insert()

# 2. Populate the newly inserted row 2 with the freshly scraped listing.
$ws.Range("A2").Value = "2025-09-16 01:14:16"
$ws.Range("B2").Value = "【急募】AI&SaaS Lineプラットフォーム開発のプロを探しています!"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5393834"
$ws.Range("G2").Value = 375
$ws.Range("H2").Value = "🔥AI,Ai ◆開発"

# 3. Every row's "取得日時" column reflects the time of this scrape run,
#    so refresh column A for all the rows that were pushed down too.
for ($r = 3; $r -le 11; $r++) {
    $ws.Cells.Item($r, 1).Value = "2025-09-16 01:14:16"
}

# 4. Column D (price) got a little wider in this revision (28 -> 30
#    characters). ColumnWidth needs a small offset subtracted to land on
#    an exact "30" once Excel re-applies its internal padding on save.
$ws.Columns.Item(4).ColumnWidth = 30 - 5/6

# 5. The row insert above does not renumber the worksheet's hyperlink
#    anchors, so rebuild the hyperlink collection from scratch, matching
#    each URL cell in F2:F11 in order (this also regenerates rId1..rId10
#    in the worksheet relationships, with the new listing becoming rId1).
$ws.Hyperlinks.Delete()
for ($r = 2; $r -le 11; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $ws.Hyperlinks.Add($cell, $cell.Value2)
}

Write-Output "Inserted new listing row; sheet now spans A1:H11."
